$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(1422, 1).Value = "Buying Opportunity"
$ws.Cells.Item(1422, 2).Value = "support Zone"
$ws.Cells.Item(1422, 3).Value = "long buildup"
$ws.Cells.Item(1422, 4).Value = "Short buildup"
$ws.Cells.Item(1422, 5).Value = "FII ENTERING"

$ws.Cells.Item(1423, 1).Value = "ACCELYA"
$ws.Cells.Item(1423, 2).Value = "ARIHANTCAP"
$ws.Cells.Item(1423, 5).Value = "LTTS"
$ws.Cells.Item(1423, 6).Value = 1852.6
$ws.Cells.Item(1423, 7).Value = 71.83
$ws.Cells.Item(1423, 10).Value = 5125.4

$ws.Cells.Item(1424, 1).Value = "ALPHAETF"
$ws.Cells.Item(1424, 2).Value = "ASIANENE"
$ws.Cells.Item(1424, 5).Value = "NATIONALUM"
$ws.Cells.Item(1424, 6).Value = 28.54
$ws.Cells.Item(1424, 7).Value = 317.05
$ws.Cells.Item(1424, 10).Value = 199.02

$ws.Cells.Item(1425, 1).Value = "APOLLOHOSP"
$ws.Cells.Item(1425, 2).Value = "BAJAJHIND"
$ws.Cells.Item(1425, 6).Value = 6328.55
$ws.Cells.Item(1425, 7).Value = 40.88

$ws.Cells.Item(1426, 1).Value = "ASMS"
$ws.Cells.Item(1426, 2).Value = "BALRAMCHIN"
$ws.Cells.Item(1426, 6).Value = 22.01
$ws.Cells.Item(1426, 7).Value = 426.6

$ws.Cells.Item(1427, 1).Value = "AURIONPRO"
$ws.Cells.Item(1427, 2).Value = "BODALCHEM"
$ws.Cells.Item(1427, 6).Value = 1609.8
$ws.Cells.Item(1427, 7).Value = 78.34999999999999

$ws.Cells.Item(1428, 1).Value = "BEML"
$ws.Cells.Item(1428, 2).Value = "HATHWAY"
$ws.Cells.Item(1428, 6).Value = 5059.95
$ws.Cells.Item(1428, 7).Value = 22.81

$ws.Cells.Item(1429, 1).Value = "BSLNIFTY"
$ws.Cells.Item(1429, 2).Value = "IEL"
$ws.Cells.Item(1429, 6).Value = 27.86
$ws.Cells.Item(1429, 7).Value = 11.8

$ws.Cells.Item(1430, 1).Value = "CANBK"
$ws.Cells.Item(1430, 2).Value = "IGPL"
$ws.Cells.Item(1430, 6).Value = 117.76
$ws.Cells.Item(1430, 7).Value = 627.75

$ws.Cells.Item(1431, 1).Value = "CENTEXT"
$ws.Cells.Item(1431, 2).Value = "JKPAPER"
$ws.Cells.Item(1431, 6).Value = 25.81
$ws.Cells.Item(1431, 7).Value = 576.15

$ws.Cells.Item(1432, 1).Value = "COASTCORP"
$ws.Cells.Item(1432, 2).Value = "KOHINOOR"
$ws.Cells.Item(1432, 6).Value = 253.02
$ws.Cells.Item(1432, 7).Value = 42.49

$ws.Cells.Item(1433, 1).Value = "GABRIEL"
$ws.Cells.Item(1433, 2).Value = "LOVABLE"
$ws.Cells.Item(1433, 6).Value = 497.3
$ws.Cells.Item(1433, 7).Value = 136.59

$ws.Cells.Item(1434, 1).Value = "GANDHITUBE"
$ws.Cells.Item(1434, 2).Value = "PAISALO"
$ws.Cells.Item(1434, 6).Value = 831.05
$ws.Cells.Item(1434, 7).Value = 73.39

$ws.Cells.Item(1435, 1).Value = "GESHIP"
$ws.Cells.Item(1435, 2).Value = "RAJMET"
$ws.Cells.Item(1435, 6).Value = 1267.8
$ws.Cells.Item(1435, 7).Value = 11.68

$ws.Cells.Item(1436, 1).Value = "GIPCL"
$ws.Cells.Item(1436, 6).Value = 252.99

$ws.Cells.Item(1437, 1).Value = "GOLDBEES"
$ws.Cells.Item(1437, 6).Value = 61.48

$ws.Cells.Item(1438, 1).Value = "GOLDETF"
$ws.Cells.Item(1438, 6).Value = 71.86

$ws.Cells.Item(1439, 1).Value = "HDFCLOWVOL"
$ws.Cells.Item(1439, 6).Value = 19.53

$ws.Cells.Item(1440, 1).Value = "HESTERBIO"
$ws.Cells.Item(1440, 6).Value = 2655.75

$ws.Cells.Item(1441, 1).Value = "HFCL"
$ws.Cells.Item(1441, 6).Value = 130.36

$ws.Cells.Item(1442, 1).Value = "IFBAGRO"
$ws.Cells.Item(1442, 6).Value = 607.8

$ws.Cells.Item(1443, 1).Value = "IRCON"
$ws.Cells.Item(1443, 6).Value = 307.75

$ws.Cells.Item(1444, 1).Value = "ITC"
$ws.Cells.Item(1444, 6).Value = 433.65

$ws.Cells.Item(1445, 1).Value = "IVZINGOLD"
$ws.Cells.Item(1445, 6).Value = 6425.4

$ws.Cells.Item(1446, 1).Value = "JMFINANCIL"
$ws.Cells.Item(1446, 6).Value = 98.23999999999999

$ws.Cells.Item(1447, 1).Value = "KELLTONTEC"
$ws.Cells.Item(1447, 6).Value = 115.81

$ws.Cells.Item(1448, 1).Value = "KOTHARIPRO"
$ws.Cells.Item(1448, 6).Value = 182.92

$ws.Cells.Item(1449, 1).Value = "LEMONTREE"
$ws.Cells.Item(1449, 6).Value = 148.15

$ws.Cells.Item(1450, 1).Value = "LTTS"
$ws.Cells.Item(1450, 6).Value = 5125.4

$ws.Cells.Item(1451, 1).Value = "LXCHEM"
$ws.Cells.Item(1451, 6).Value = 268.45

$ws.Cells.Item(1452, 1).Value = "MAHKTECH"
$ws.Cells.Item(1452, 6).Value = 14.03

$ws.Cells.Item(1453, 1).Value = "MANGALAM"
$ws.Cells.Item(1453, 6).Value = 113.32

$ws.Cells.Item(1454, 1).Value = "MARICO"
$ws.Cells.Item(1454, 6).Value = 615.35

$ws.Cells.Item(1455, 1).Value = "MOM100"
$ws.Cells.Item(1455, 6).Value = 61.07

$ws.Cells.Item(1456, 1).Value = "MONIFTY500"
$ws.Cells.Item(1456, 6).Value = 23.62

$ws.Cells.Item(1457, 1).Value = "NATIONALUM"
$ws.Cells.Item(1457, 6).Value = 199.02

$ws.Cells.Item(1458, 1).Value = "NIACL"
$ws.Cells.Item(1458, 6).Value = 273.84

$ws.Cells.Item(1459, 1).Value = "NOCIL"
$ws.Cells.Item(1459, 6).Value = 306.4

$ws.Cells.Item(1460, 1).Value = "NURECA"
$ws.Cells.Item(1460, 6).Value = 282.4

$ws.Cells.Item(1461, 1).Value = "PLASTIBLEN"
$ws.Cells.Item(1461, 6).Value = 299.3

$ws.Cells.Item(1462, 1).Value = "RKEC"
$ws.Cells.Item(1462, 6).Value = 132.29

$ws.Cells.Item(1463, 1).Value = "RUBYMILLS"
$ws.Cells.Item(1463, 6).Value = 245.27

$ws.Cells.Item(1464, 1).Value = "RVNL"
$ws.Cells.Item(1464, 6).Value = 491.05

$ws.Cells.Item(1465, 1).NumberFormat = "@"
$ws.Cells.Item(1465, 1).Value = "05/07/2024"

$ws.Cells.Item(1466, 1).Value = "Buying Opportunity"
$ws.Cells.Item(1466, 2).Value = "support Zone"
$ws.Cells.Item(1466, 3).Value = "long buildup"
$ws.Cells.Item(1466, 4).Value = "Short buildup"
$ws.Cells.Item(1466, 5).Value = "FII ENTERING"

$ws.Cells.Item(1467, 1).Value = "AJMERA"
$ws.Cells.Item(1467, 2).Value = "ADSL"
$ws.Cells.Item(1467, 3).Value = "BHEL"
$ws.Cells.Item(1467, 5).Value = "BEL"
$ws.Cells.Item(1467, 6).Value = 746.4
$ws.Cells.Item(1467, 7).Value = 218.09
$ws.Cells.Item(1467, 8).Value = 328.35
$ws.Cells.Item(1467, 10).Value = 334.6

$ws.Cells.Item(1468, 1).Value = "ALKEM"
$ws.Cells.Item(1468, 2).Value = "ATL"
$ws.Cells.Item(1468, 3).Value = "NATIONALUM"
$ws.Cells.Item(1468, 5).Value = "IRCTC"
$ws.Cells.Item(1468, 6).Value = 5221.85
$ws.Cells.Item(1468, 7).Value = 50.61
$ws.Cells.Item(1468, 8).Value = 206.31
$ws.Cells.Item(1468, 10).Value = 1045.6

$ws.Cells.Item(1469, 1).Value = "ALPHAETF"
$ws.Cells.Item(1469, 2).Value = "CHEMPLASTS"
$ws.Cells.Item(1469, 5).Value = "ITC"
$ws.Cells.Item(1469, 6).Value = 28.71
$ws.Cells.Item(1469, 7).Value = 543.9
$ws.Cells.Item(1469, 10).Value = 443.6

$ws.Cells.Item(1470, 1).Value = "ANANDRATHI"
$ws.Cells.Item(1470, 2).Value = "DELHIVERY"
$ws.Cells.Item(1470, 5).Value = "NATIONALUM"
$ws.Cells.Item(1470, 6).Value = 4094.4
$ws.Cells.Item(1470, 7).Value = 394.45
$ws.Cells.Item(1470, 10).Value = 206.31

$ws.Cells.Item(1471, 1).Value = "ANANTRAJ"
$ws.Cells.Item(1471, 2).Value = "DLINKINDIA"
$ws.Cells.Item(1471, 6).Value = 544.3
$ws.Cells.Item(1471, 7).Value = 531.35

$ws.Cells.Item(1472, 1).Value = "BEL"
$ws.Cells.Item(1472, 2).Value = "ENIL"
$ws.Cells.Item(1472, 6).Value = 334.6
$ws.Cells.Item(1472, 7).Value = 248.58

$ws.Cells.Item(1473, 1).Value = "GESHIP"
$ws.Cells.Item(1473, 2).Value = "GPPL"
$ws.Cells.Item(1473, 6).Value = 1343.2
$ws.Cells.Item(1473, 7).Value = 222.38

$ws.Cells.Item(1474, 1).Value = "HBLPOWER"
$ws.Cells.Item(1474, 2).Value = "JKPAPER"
$ws.Cells.Item(1474, 6).Value = 597.85
$ws.Cells.Item(1474, 7).Value = 569.9

$ws.Cells.Item(1475, 1).Value = "HDFCBSE500"
$ws.Cells.Item(1475, 2).Value = "KAKATCEM"
$ws.Cells.Item(1475, 6).Value = 36.91
$ws.Cells.Item(1475, 7).Value = 228.11

$ws.Cells.Item(1476, 1).Value = "HDFCSML250"
$ws.Cells.Item(1476, 2).Value = "KINGFA"
$ws.Cells.Item(1476, 6).Value = 178.32
$ws.Cells.Item(1476, 7).Value = 2500.5

$ws.Cells.Item(1477, 1).Value = "HINDCOPPER"
$ws.Cells.Item(1477, 2).Value = "LOVABLE"
$ws.Cells.Item(1477, 6).Value = 342.25
$ws.Cells.Item(1477, 7).Value = 134.09

$ws.Cells.Item(1478, 1).Value = "HPAL"
$ws.Cells.Item(1478, 2).Value = "LTIM"
$ws.Cells.Item(1478, 6).Value = 101.39
$ws.Cells.Item(1478, 7).Value = 5389.7

$ws.Cells.Item(1479, 1).Value = "IRCTC"
$ws.Cells.Item(1479, 2).Value = "MARALOVER"
$ws.Cells.Item(1479, 6).Value = 1045.6
$ws.Cells.Item(1479, 7).Value = 72.61

$ws.Cells.Item(1480, 1).Value = "ITC"
$ws.Cells.Item(1480, 2).Value = "NUVAMA"
$ws.Cells.Item(1480, 6).Value = 443.6
$ws.Cells.Item(1480, 7).Value = 4715.45

$ws.Cells.Item(1481, 1).Value = "ITI"
$ws.Cells.Item(1481, 2).Value = "POKARNA"
$ws.Cells.Item(1481, 6).Value = 312.7
$ws.Cells.Item(1481, 7).Value = 729.5

$ws.Cells.Item(1482, 1).Value = "IVZINGOLD"
$ws.Cells.Item(1482, 6).Value = 6465

$ws.Cells.Item(1483, 1).Value = "JBMA"
$ws.Cells.Item(1483, 6).Value = 2231.8

$ws.Cells.Item(1484, 1).Value = "JITFINFRA"
$ws.Cells.Item(1484, 6).Value = 831.75

$ws.Cells.Item(1485, 1).Value = "KELLTONTEC"
$ws.Cells.Item(1485, 6).Value = 122.14

$ws.Cells.Item(1486, 1).Value = "KIRIINDUS"
$ws.Cells.Item(1486, 6).Value = 363.85

$ws.Cells.Item(1487, 1).Value = "KSCL"
$ws.Cells.Item(1487, 6).Value = 980.7

$ws.Cells.Item(1488, 1).Value = "LEXUS"
$ws.Cells.Item(1488, 6).Value = 48.93

$ws.Cells.Item(1489, 1).Value = "LPDC"
$ws.Cells.Item(1489, 6).Value = 15.24

$ws.Cells.Item(1490, 1).Value = "MIDHANI"
$ws.Cells.Item(1490, 6).Value = 525.95

$ws.Cells.Item(1491, 1).Value = "MIRCELECTR"
$ws.Cells.Item(1491, 6).Value = 21.29

$ws.Cells.Item(1492, 1).Value = "MON100"
$ws.Cells.Item(1492, 6).Value = 167.53

$ws.Cells.Item(1493, 1).Value = "MONQ50"
$ws.Cells.Item(1493, 6).Value = 68.41

$ws.Cells.Item(1494, 1).Value = "NATIONALUM"
$ws.Cells.Item(1494, 6).Value = 206.31

$ws.Cells.Item(1495, 1).Value = "NETWEB"
$ws.Cells.Item(1495, 6).Value = 2685.65

$ws.Cells.Item(1496, 1).Value = "NV20BEES"
$ws.Cells.Item(1496, 6).Value = 153.98

$ws.Cells.Item(1497, 1).Value = "PAGEIND"
$ws.Cells.Item(1497, 6).Value = 39460.25

$ws.Cells.Item(1498, 1).Value = "PNC"
$ws.Cells.Item(1498, 6).Value = 68.04000000000001

$ws.Cells.Item(1499, 1).Value = "RADIANTCMS"
$ws.Cells.Item(1499, 6).Value = 80.16

$ws.Cells.Item(1500, 1).NumberFormat = "@"
$ws.Cells.Item(1500, 1).Value = "08/07/2024"
